$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "6.68") must be forced
# to stay text, matching the source inlineStr cells. We flip NumberFormat to
# text ("@") before assigning, then restore the cell style to Normal so the
# on-disk style index is unchanged (avoids introducing spurious style diffs).
$textForcedCells = @("D5", "D6", "D10", "D12", "D16", "D18", "D19", "D22", "D24", "D28", "D32", "D33", "D34", "D37", "D39", "D40", "D44", "D45", "D51")
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "59.755.99"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "2.649.92"
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "537.14"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "145.62"
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("D9").Value = "2.665.88"
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("D10").Value = "6.68"
$ws.Range("E10").Value = "  +3.40%  "
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").Value = "0.339"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("D14").Value = "3.116.43"
$ws.Range("E14").Value = "  +1.94%  "
$ws.Range("D15").Value = "59.673.37"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("D16").Value = "21.18"
$ws.Range("E16").Value = "  +2.87%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.635.34"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.0000135"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").Value = "344.90"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("E20").Value = "  +1.79%  "
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("D22").Value = "6.35"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "66.65"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("E25").Value = "  +2.17%  "
$ws.Range("E26").Value = "  -1.57%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Value = "7.29"
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("D29").Value = "0.0₃0753"
$ws.Range("E29").Value = "  +2.03%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  +1.42%  "
$ws.Range("D32").Value = "5.84"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").Value = "19.00"
$ws.Range("E33").Value = "  +0.96%  "
$ws.Range("D34").Value = "149.95"
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("E35").Value = "  +0.82%  "
$ws.Range("E36").Value = "  +2.21%  "
$ws.Range("D37").Value = "0.843"
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("D39").Value = "0.826"
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("D40").Value = "292.50"
$ws.Range("E40").Value = "  +5.66%  "
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("E43").Value = "  +1.30%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "19.49"
$ws.Range("E44").Value = "  +5.27%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "0.0543"
$ws.Range("E45").Value = "  +4.20%  "
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("D48").Value = "1.974.35"
$ws.Range("E48").Value = "  +1.60%  "
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("D51").Value = "18.37"
$ws.Range("E51").Value = "  +0.17%  "

foreach ($addr in $textForcedCells) {
    $ws.Range($addr).Style = "Normal"
}
